$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'70.015.16"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = "'3.785.08"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'616.88"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.43%  '
$ws.Range('D6').Value = "'178.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.50%  '
$ws.Range('D7').Value = "'3.787.50"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.36%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('D11').Value = "'6.33"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.47%  '
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').Value = "'40.91"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').Value = "'0.0000256"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').Value = "'4.411.57"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.91%  '
$ws.Range('D16').Value = "'3.783.03"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').Value = "'70.085.03"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = "'7.59"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = "'514.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = "'16.62"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.60%  '
$ws.Range('D22').Value = "'9.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.18%  '
$ws.Range('E23').Value = '  -3.01%  '
$ws.Range('D24').Value = "'2.52"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.17%  '
$ws.Range('D25').Value = "'88.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = "'13.31"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').Value = "'11.10"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('D28').Value = "'0.0000135"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +24.42%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('E31').Value = '  -5.11%  '
$ws.Range('D32').Value = "'2.84"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('D33').Value = "'31.81"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('E34').Value = '  -2.00%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = "'6.22"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').Value = "'0.340"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('D41').Value = "'51.32"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('D42').Value = "'44.41"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.48%  '
$ws.Range('D43').Value = "'8.78"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').Value = "'423.74"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.64%  '
$ws.Range('D45').Value = "'3.065.60"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('D46').Value = "'2.75"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = "'27.67"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('D51').Value = "'135.25"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.07%  '
